# Fill A1:A35 with the numbers 1..35 and select that range,
# matching the new sheet1.xml data / dimension / selection in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 35; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}

$ws.Range("A1:A35").Select()
